$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "discussion summary" cell text to mention the formal rotation
# of the group leader, and refresh the adjacent duration cell so the
# shared-string table settles in the same order as the authored workbook.
$ws.Range("C7").Value = "集体检查产物，最后整理，阶段总结，正式轮换组长。"
$ws.Range("D7").Value = "60分钟"

# Move the sheet's active selection from D10 to F4, matching the author's
# last cursor position when the workbook was saved.
$ws.Range("F4").Select() | Out-Null
